$d = $word.ActiveDocument

# The letter heading in the first table cell currently reads
# "*Karl an Ferdinand." -- strip the leading asterisk so it reads
# "Karl an Ferdinand." (the rest of the run's text/formatting stays
# untouched, so only the single stray character is removed).
$cell = $d.Tables(1).Cell(1, 1).Range
$star = $d.Range($cell.Start, $cell.Start + 1)

if ($star.Text -eq "*") {
    $star.Delete()
} else {
    # Fallback in case the cell layout ever changes: a plain literal
    # find & replace of the whole heading text.
    $d.Content.Find.Execute("*Karl an Ferdinand.", $false, $false, $false, $false, $false, `
                             $true, 1, $false, "Karl an Ferdinand.", 2)
}
